# "Elimina EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The detail table (rows 16-51, columns C:G) is rebuilt: instead of being
# interleaved by period across the three workers, the rows are now grouped
# by worker (document id / name) with each worker's periods listed in
# descending order (most recent period first).
#
# Columns: C = Tipo/N Doc Trabajador (doc number), D = Nombre Trabajador,
#          E = Periodo Mora, F = Valor Mora, G = Salario Basico

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Each worker's block: document number, name, and list of (period, valorMora, salarioBasico)
$workers = @(
    @{
        Doc = "45689110"
        Nombre = "INGRID DEL CARMEN PUELLO PEREZ"
        Periodos = @(
            @("2211", 25439, 908526),
            @("2210", 36341, 908526),
            @("2209", 36341, 908526),
            @("2208", 36341, 908526),
            @("2207", 36341, 908526),
            @("2206", 36341, 908526),
            @("2205", 36341, 908526),
            @("2204", 36341, 908526),
            @("2203", 36341, 908526),
            @("2202", 36341, 908526),
            @("2201", 36341, 908526),
            @("2112", 36341, 908526),
            @("2111", 36341, 908526),
            @("2110", 36341, 908526),
            @("2109", 36341, 908526),
            @("2108", 36341, 908526)
        )
    },
    @{
        Doc = "45551768"
        Nombre = "KAROL PUELLO PEREZ"
        Periodos = @(
            @("2211", 25439, 1000000),
            @("2210", 40000, 1000000),
            @("2209", 40000, 1000000),
            @("2208", 40000, 1000000),
            @("2207", 40000, 1000000),
            @("2206", 40000, 1000000),
            @("2205", 40000, 1000000)
        )
    },
    @{
        Doc = "1047490024"
        Nombre = "LUIS DAVID MARTINEZ GASPAR"
        Periodos = @(
            @("2211", 25439, 908526),
            @("2210", 36341, 908526),
            @("2209", 36341, 908526),
            @("2208", 36341, 908526),
            @("2207", 36341, 908526),
            @("2206", 36341, 908526),
            @("2205", 36341, 908526),
            @("2204", 36341, 908526),
            @("2203", 36341, 908526),
            @("2202", 36341, 908526),
            @("2201", 36341, 908526),
            @("2112", 36341, 908526),
            @("2111", 36341, 908526)
        )
    }
)

$row = 16
foreach ($worker in $workers) {
    foreach ($periodo in $worker.Periodos) {
        $ws.Range("C$row").Value = $worker.Doc
        $ws.Range("D$row").Value = $worker.Nombre
        $ws.Range("E$row").Value = $periodo[0]
        $ws.Range("F$row").Value = $periodo[1]
        $ws.Range("G$row").Value = $periodo[2]
        $row = $row + 1
    }
}
